# "covid daily emissions update"
# Applies new rows/formulas to the gas_trend worksheet of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("gas_trend")
$ws.Activate()

# --- Row 8: new AH8 = AE5-AD5 (style matches AG8 -> numberformat 0.0) ---
$ws.Range("AH8").Formula = "=AE5-AD5"
$ws.Range("AH8").NumberFormat = "0.0"

# --- Row 11: fill V11:AD11 with shared SUM formula, AE11 with its own SUM formula ---
$ws.Range("V11:AD11").Formula = "=SUM(V5:V6)"
$ws.Range("V11:AD11").NumberFormat = "0"
$ws.Range("AE11").Formula = "=SUM(AE5:AE6)"

# AF11: CAGR-like growth calc between AE11 and V11
$ws.Range("AF11").Formula = "=((AE11/V11)^(AE1-V1))-1"
$ws.Range("AF11").NumberFormat = "0.00"

# --- Row 12 (new row): AF12 similarly vs AD11 ---
$ws.Range("AF12").Formula = "=((AD11/V11)^(AD1-V1))-1"
$ws.Range("AF12").NumberFormat = "0.00"

# --- Rows 14-18: convert the per-cell formulas into shared formulas (values unchanged) ---
$ws.Range("B14:AE14").Formula = '=(B2*$AH2)^2'
$ws.Range("B15:AE15").Formula = '=(B3*$AH3)^2'
$ws.Range("B16:AE16").Formula = '=(B4*$AH4)^2'
$ws.Range("B17:AE17").Formula = '=(B5*$AH5)^2'
$ws.Range("B18:AE18").Formula = '=(B6*$AH6)^2'

# --- Row 18: new AF18 cell with SQRT(SUM(V17:V18)) ---
$ws.Range("AF18").Formula = "=SQRT(SUM(V17:V18))"

# --- New rows 21-30: growth-since-1990 summary block ---
$ws.Range("AE21").Value = "abs growth since 1990"

$ws.Range("AE22").Formula = "=(AE2-B2)/B2"
$ws.Range("AE22").NumberFormat = "0%"

# AE23:AE28 share one formula group (relative refs shift row by row), but row 27
# is intentionally left blank, so fill the whole block then clear just that cell.
$ws.Range("AE23:AE28").Formula = "=(AE3-B3)/B3"
$ws.Range("AE23:AE28").NumberFormat = "0%"
$ws.Range("AE27").ClearContents()

$ws.Range("AE30").Value = "share of gas in total"

# --- Restore the active selection to match the final view state ---
$ws.Range("AH8").Select()
